$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking") updates
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 ("Total") updates
$ws.Range("B12").Value = 85
$ws.Range("C12").Value = -1.2
$ws.Range("E12").Value = "83.8/140"
